$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.440334666666667
$ws.Range("H2").Value = 10.321004
$ws.Range("I2").Value = 0.03090290794544385
$ws.Range("J2").Value = 0.03090290794544386
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 120.5730638824218
$ws.Range("R2").Value = 1085.157574941796
$ws.Range("S2").Value = 0.01181921855857787
$ws.Range("T2").Value = 0.01181921855857788
$ws.Range("G3").Value = 3.440334666666667
$ws.Range("H3").Value = 10.321004
$ws.Range("I3").Value = 0.03090290794544385
$ws.Range("J3").Value = 0.03090290794544386
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 102.9116563340253
$ws.Range("R3").Value = 926.204907006228
$ws.Range("S3").Value = 0.01008795264274966
$ws.Range("T3").Value = 0.01008795264274966
$ws.Range("G4").Value = 3.440334666666667
$ws.Range("H4").Value = 10.321004
$ws.Range("I4").Value = 0.03090290794544385
$ws.Range("J4").Value = 0.03090290794544386
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 30.65932563152577
$ws.Range("R4").Value = 275.933930683732
$ws.Range("S4").Value = 0.003005391576106752
$ws.Range("T4").Value = 0.003005391576106753
$ws.Range("G5").Value = 3.440334666666667
$ws.Range("H5").Value = 10.321004
$ws.Range("I5").Value = 0.03090290794544385
$ws.Range("J5").Value = 0.03090290794544386
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 61.11015436769111
$ws.Range("R5").Value = 549.99138930922
$ws.Range("S5").Value = 0.005990345168009561
$ws.Range("T5").Value = 0.005990345168009562
$ws.Range("I6").Value = 0.73221566931385
$ws.Range("J6").Value = 0.7322156693138502
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 2856.866636232056
$ws.Range("R6").Value = 25711.79972608851
$ws.Range("S6").Value = 0.2800453938805362
$ws.Range("T6").Value = 0.2800453938805363
$ws.Range("I7").Value = 0.73221566931385
$ws.Range("J7").Value = 0.7322156693138502
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.2390246577881159
$ws.Range("T7").Value = 0.2390246577881159
$ws.Range("I8").Value = 0.73221566931385
$ws.Range("J8").Value = 0.7322156693138502
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 726.4442128757242
$ws.Range("R8").Value = 6537.997915881519
$ws.Range("S8").Value = 0.07120995889235257
$ws.Range("T8").Value = 0.07120995889235258
$ws.Range("I9").Value = 0.73221566931385
$ws.Range("J9").Value = 0.7322156693138502
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 1447.948285682568
$ws.Range("R9").Value = 13031.53457114312
$ws.Range("S9").Value = 0.1419356587528453
$ws.Range("T9").Value = 0.1419356587528453
$ws.Range("G10").Value = 25.13705366666667
$ws.Range("H10").Value = 75.411161
$ws.Range("I10").Value = 0.2257943283853049
$ws.Range("J10").Value = 0.225794328385305
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 880.9757977712823
$ws.Range("R10").Value = 7928.782179941541
$ws.Range("S10").Value = 0.08635797385749527
$ws.Range("T10").Value = 0.08635797385749529
$ws.Range("G11").Value = 25.13705366666667
$ws.Range("H11").Value = 75.411161
$ws.Range("I11").Value = 0.2257943283853049
$ws.Range("J11").Value = 0.225794328385305
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 751.9314481984364
$ws.Range("R11").Value = 6767.383033785927
$ws.Range("S11").Value = 0.07370835442974058
$ws.Range("T11").Value = 0.07370835442974058
$ws.Range("G12").Value = 25.13705366666667
$ws.Range("H12").Value = 75.411161
$ws.Range("I12").Value = 0.2257943283853049
$ws.Range("J12").Value = 0.225794328385305
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 224.0145766197181
$ws.Range("R12").Value = 2016.131189577463
$ws.Range("S12").Value = 0.02195911056849024
$ws.Range("T12").Value = 0.02195911056849024
$ws.Range("G13").Value = 25.13705366666667
$ws.Range("H13").Value = 75.411161
$ws.Range("I13").Value = 0.2257943283853049
$ws.Range("J13").Value = 0.225794328385305
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 446.5057556180395
$ws.Range("R13").Value = 4018.551800562356
$ws.Range("S13").Value = 0.04376888952957882
$ws.Range("T13").Value = 0.04376888952957883
$ws.Range("G14").Value = 1.234295333333334
$ws.Range("H14").Value = 3.702886
$ws.Range("I14").Value = 0.01108709435540116
$ws.Range("J14").Value = 0.01108709435540116
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 43.25822470636823
$ws.Range("R14").Value = 389.3240223573141
$ws.Range("S14").Value = 0.00424040325258068
$ws.Range("T14").Value = 0.00424040325258068
$ws.Range("G15").Value = 1.234295333333334
$ws.Range("H15").Value = 3.702886
$ws.Range("I15").Value = 0.01108709435540116
$ws.Range("J15").Value = 0.01108709435540116
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 36.92180833144467
$ws.Range("R15").Value = 332.296274983002
$ws.Range("S15").Value = 0.003619273726616202
$ws.Range("T15").Value = 0.003619273726616202
$ws.Range("G16").Value = 1.234295333333334
$ws.Range("H16").Value = 3.702886
$ws.Range("I16").Value = 0.01108709435540116
$ws.Range("J16").Value = 0.01108709435540116
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 10.99970387090422
$ws.Range("R16").Value = 98.997334838138
$ws.Range("S16").Value = 0.001078249983401191
$ws.Range("T16").Value = 0.001078249983401192
$ws.Range("G17").Value = 1.234295333333334
$ws.Range("H17").Value = 3.702886
$ws.Range("I17").Value = 0.01108709435540116
$ws.Range("J17").Value = 0.01108709435540116
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 21.92460491885889
$ws.Range("R17").Value = 197.32144426973
$ws.Range("S17").Value = 0.00214916739280309
$ws.Range("T17").Value = 0.00214916739280309
